<#
  scene_cat_exp memory_kitchens sheet cleanup.

  The stimuli/ratings columns (H..V) of rows 2-25 (plus two rows further
  down, 32 and 41, which swap stimuli with rows 7 and 19) get shuffled:
  each rows old catch/target/new image+rating data is replaced by the
  data that used to live a bit further down the list, the old catch trial
  img (stimuli/catch_05_supermarket.jpg, row 25) is dropped, and a new
  catch trial (stimuli/catch_28.jpg) is inserted at row 21.
#>

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 'stimuli/img_inqod.png'
$ws.Range("M2").Value = 70.84848484848484
$ws.Range("N2").Value = 50.63636363636363
$ws.Range("O2").Value = 60.74242424242424
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 5

# Row 3
$ws.Range("I3").Value = 'target'
$ws.Range("J3").Value = 'old'
$ws.Range("K3").Value = 'j'
$ws.Range("L3").Value = 'stimuli/img_aplao.png'
$ws.Range("M3").Value = 64.0909090909091
$ws.Range("N3").Value = 40.75757575757576
$ws.Range("O3").Value = 52.42424242424242
$ws.Range("P3").Value = 33
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 3
$ws.Range("V3").Value = 3

# Row 4
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = 'new'
$ws.Range("K4").Value = 'f'
$ws.Range("L4").Value = 'stimuli/img_mgnmm.png'
$ws.Range("M4").Value = 79.1470588235294
$ws.Range("N4").Value = 60.38235294117647
$ws.Range("O4").Value = 69.76470588235294
$ws.Range("P4").Value = 34
$ws.Range("Q4").Value = 8
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 8
$ws.Range("V4").Value = 8

# Row 5
$ws.Range("I5").Value = 'target'
$ws.Range("J5").Value = 'old'
$ws.Range("K5").Value = 'j'
$ws.Range("L5").Value = 'stimuli/img_30vz5.png'
$ws.Range("M5").Value = 86.21212121212122
$ws.Range("N5").Value = 68.27272727272727
$ws.Range("O5").Value = 77.24242424242425
$ws.Range("P5").Value = 33
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 10
$ws.Range("V5").Value = 10

# Row 6
$ws.Range("L6").Value = 'stimuli/img_esb4r.png'
$ws.Range("M6").Value = 60.73529411764706
$ws.Range("N6").Value = 38.58823529411764
$ws.Range("O6").Value = 49.66176470588235
$ws.Range("P6").Value = 34
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 3
$ws.Range("V6").Value = 3

# Row 7
$ws.Range("L7").Value = 'stimuli/img_mjxmq.png'
$ws.Range("M7").Value = 77.07692307692308
$ws.Range("N7").Value = 58.15384615384615
$ws.Range("O7").Value = 67.61538461538461
$ws.Range("P7").Value = 39
$ws.Range("Q7").Value = 7
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 7
$ws.Range("T7").Value = 7
$ws.Range("U7").Value = 7
$ws.Range("V7").Value = 7

# Row 8
$ws.Range("L8").Value = 'stimuli/img_qmgwq.png'
$ws.Range("M8").Value = 84.58333333333333
$ws.Range("N8").Value = 64.44444444444444
$ws.Range("O8").Value = 74.51388888888889
$ws.Range("P8").Value = 36
$ws.Range("Q8").Value = 9
$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 9
$ws.Range("T8").Value = 9
$ws.Range("U8").Value = 9
$ws.Range("V8").Value = 9

# Row 9
$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = 'new'
$ws.Range("K9").Value = 'f'
$ws.Range("L9").Value = 'stimuli/img_z293c.png'
$ws.Range("M9").Value = 71.26470588235294
$ws.Range("N9").Value = 46.88235294117647
$ws.Range("O9").Value = 59.07352941176471
$ws.Range("P9").Value = 34
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 5

# Row 10
$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = 'new'
$ws.Range("K10").Value = 'f'
$ws.Range("L10").Value = 'stimuli/img_vbrb7.png'
$ws.Range("M10").Value = 85.5625
$ws.Range("N10").Value = 71.46875
$ws.Range("O10").Value = 78.515625
$ws.Range("P10").Value = 32
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = 10
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 10
$ws.Range("V10").Value = 10

# Row 11
$ws.Range("L11").Value = 'stimuli/img_iyxnj.png'
$ws.Range("M11").Value = 75.30555555555556
$ws.Range("N11").Value = 54.33333333333334
$ws.Range("O11").Value = 64.81944444444444
$ws.Range("P11").Value = 36
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 6
$ws.Range("T11").Value = 6
$ws.Range("U11").Value = 6
$ws.Range("V11").Value = 6

# Row 12
$ws.Range("L12").Value = 'stimuli/img_ce9vx.png'
$ws.Range("M12").Value = 75.9090909090909
$ws.Range("N12").Value = 57.12121212121212
$ws.Range("O12").Value = 66.51515151515152
$ws.Range("P12").Value = 33
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7

# Row 13
$ws.Range("I13").Value = 'target'
$ws.Range("J13").Value = 'old'
$ws.Range("K13").Value = 'j'
$ws.Range("L13").Value = 'stimuli/img_cv6mf.png'
$ws.Range("M13").Value = 66.8
$ws.Range("N13").Value = 42.08
$ws.Range("O13").Value = 54.44
$ws.Range("P13").Value = 25

# Row 14
$ws.Range("I14").Value = 'target'
$ws.Range("J14").Value = 'old'
$ws.Range("K14").Value = 'j'
$ws.Range("L14").Value = 'stimuli/img_p3hpc.png'
$ws.Range("M14").Value = 72.83333333333333
$ws.Range("N14").Value = 52.22222222222222
$ws.Range("O14").Value = 62.52777777777777
$ws.Range("P14").Value = 36
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 6
$ws.Range("S14").Value = 6
$ws.Range("T14").Value = 6
$ws.Range("U14").Value = 6
$ws.Range("V14").Value = 6

# Row 15
$ws.Range("L15").Value = 'stimuli/img_lszzj.png'
$ws.Range("M15").Value = 64.70588235294117
$ws.Range("N15").Value = 45.58823529411764
$ws.Range("O15").Value = 55.14705882352941
$ws.Range("P15").Value = 34
$ws.Range("Q15").Value = 4
$ws.Range("R15").Value = 4
$ws.Range("S15").Value = 4
$ws.Range("T15").Value = 4
$ws.Range("U15").Value = 4
$ws.Range("V15").Value = 4

# Row 16
$ws.Range("L16").Value = 'stimuli/img_7ed9m.png'
$ws.Range("M16").Value = 80.71875
$ws.Range("N16").Value = 58.65625
$ws.Range("O16").Value = 69.6875
$ws.Range("Q16").Value = 8
$ws.Range("R16").Value = 8
$ws.Range("S16").Value = 8
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = 8
$ws.Range("V16").Value = 8

# Row 17
$ws.Range("L17").Value = 'stimuli/img_r2lxk.png'
$ws.Range("M17").Value = 89.24242424242425
$ws.Range("N17").Value = 67.6969696969697
$ws.Range("O17").Value = 78.46969696969697
$ws.Range("P17").Value = 33
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = 10
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = 10
$ws.Range("U17").Value = 10
$ws.Range("V17").Value = 10

# Row 19
$ws.Range("L19").Value = 'stimuli/img_ewrjk.png'
$ws.Range("M19").Value = 73.0909090909091
$ws.Range("N19").Value = 53.39393939393939
$ws.Range("O19").Value = 63.24242424242424
$ws.Range("P19").Value = 33
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = 6
$ws.Range("S19").Value = 6
$ws.Range("T19").Value = 6
$ws.Range("U19").Value = 6
$ws.Range("V19").Value = 6

# Row 21
$ws.Range("H21").ClearContents()
$ws.Range("J21").Value = 'catch'
$ws.Range("L21").Value = 'stimuli/catch_28.jpg'
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("O21").ClearContents()
$ws.Range("P21").ClearContents()
$ws.Range("Q21").ClearContents()
$ws.Range("R21").ClearContents()
$ws.Range("S21").ClearContents()
$ws.Range("T21").ClearContents()
$ws.Range("U21").ClearContents()

# Row 22
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = 'new'
$ws.Range("K22").Value = 'f'
$ws.Range("L22").Value = 'stimuli/img_kwxq1.png'
$ws.Range("M22").Value = 68.53125
$ws.Range("N22").Value = 44.09375
$ws.Range("O22").Value = 56.3125
$ws.Range("P22").Value = 32
$ws.Range("Q22").Value = 4
$ws.Range("R22").Value = 4
$ws.Range("S22").Value = 4
$ws.Range("T22").Value = 4
$ws.Range("U22").Value = 4
$ws.Range("V22").Value = 4

# Row 23
$ws.Range("L23").Value = 'stimuli/img_wppku.png'
$ws.Range("M23").Value = 75.02941176470588
$ws.Range("N23").Value = 53.05882352941177
$ws.Range("O23").Value = 64.04411764705883
$ws.Range("P23").Value = 34

# Row 24
$ws.Range("I24").ClearContents()
$ws.Range("J24").Value = 'new'
$ws.Range("K24").Value = 'f'
$ws.Range("L24").Value = 'stimuli/img_7ucnr.png'
$ws.Range("M24").Value = 70.39393939393939
$ws.Range("N24").Value = 47.90909090909091
$ws.Range("O24").Value = 59.15151515151515
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 5

# Row 25
$ws.Range("H25").Value = 'kitchens'
$ws.Range("I25").Value = 'target'
$ws.Range("J25").Value = 'old'
$ws.Range("K25").Value = 'j'
$ws.Range("L25").Value = 'stimuli/img_7wul8.png'
$ws.Range("M25").Value = 43.03030303030303
$ws.Range("N25").Value = 25.54545454545455
$ws.Range("O25").Value = 34.28787878787879
$ws.Range("P25").Value = 33
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 1
$ws.Range("S25").Value = 1
$ws.Range("T25").Value = 1
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = 1

# Row 32
$ws.Range("L32").Value = 'stimuli/img_mawe6.png'
$ws.Range("M32").Value = 83.48387096774194
$ws.Range("N32").Value = 65.54838709677419
$ws.Range("O32").Value = 74.51612903225806
$ws.Range("P32").Value = 31
$ws.Range("Q32").Value = 9
$ws.Range("R32").Value = 9
$ws.Range("S32").Value = 9
$ws.Range("T32").Value = 9
$ws.Range("U32").Value = 9
$ws.Range("V32").Value = 9

# Row 41
$ws.Range("L41").Value = 'stimuli/img_zi8qc.png'
$ws.Range("M41").Value = 77.14285714285714
$ws.Range("N41").Value = 57.02857142857143
$ws.Range("O41").Value = 67.08571428571429
$ws.Range("P41").Value = 35
$ws.Range("Q41").Value = 7
$ws.Range("R41").Value = 7
$ws.Range("S41").Value = 7
$ws.Range("T41").Value = 7
$ws.Range("U41").Value = 7
$ws.Range("V41").Value = 7
